$wb = $excel.ActiveWorkbook

# --- createcourse sheet: insert a new "CREDITS_VL" column between
#     COURSE_NAME_VL and DESCRIPTION_VL, with value 3 for the sample row ---
$ws2 = $wb.Worksheets.Item("createcourse")
$ws2.Columns("D:D").Insert()
$ws2.Range("D1").Value = "CREDITS_VL"
$ws2.Range("D2").Value = 3
[void]$ws2.Range("F4").Select()

# --- clo sheet: fix the "ESCRIPTION_VL" header typo -> "DESCRIPTION_VL" ---
$ws3 = $wb.Worksheets.Item("clo")
$ws3.Range("C1").Value = "DESCRIPTION_VL"

# clo becomes the active / selected sheet and cell
[void]$ws3.Activate()
[void]$ws3.Range("I3").Select()
